$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.238.62'
$ws.Range("E2").Value = '  -4.76%  '
$ws.Range("D3").Value = '1.743.17'
$ws.Range("E3").Value = '  -4.93%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9945'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.64%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9990'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4310'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3588'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07052'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8243'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.54%  '
$ws.Range("D12").Value = '1.745.91'
$ws.Range("E12").Value = '  -10.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.159'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.242'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06757'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008580'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9991'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.31%  '
$ws.Range("D21").Value = '26.379.49'
$ws.Range("E21").Value = '  -5.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.935'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("D24").Value = '1.979.79'
$ws.Range("E24").Value = '  -8.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.880'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.976'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.625'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08877'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7050'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.246'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.763'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.36%  '
$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.000'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.070'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.059'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.59%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01861'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.89%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05016'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4836'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1579'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.486'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.091'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.784'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.50%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9956'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.965'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.65%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06158'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4411'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.546'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.686'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.04%  '
